# Generate Report for Handoff
# Updates the "b.md" rows across the Overview, zh-cn and de-de sheets to
# reflect that the file is now ready for handoff (instead of "Handed back"),
# with a freshly generated handoff file name + timestamp.

$wb = $excel.ActiveWorkbook

$newStatus = "Ready for handoff"

# ---------------------------------------------------------------------
# Overview sheet: row 3 is the "b.md" row.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus
$wsOverview.Range("D3").Value = "2016-03-22 18:33:50"

# ---------------------------------------------------------------------
# zh-cn sheet: row 3 is the "b.md" row.
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $newStatus
$wsZhCn.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("E3").Value = "2016-03-22 18:33:45"

foreach ($hl in $wsZhCn.Hyperlinks) {
    if ($hl.Range.Address() -eq "`$D`$3") {
        $hl.TextToDisplay = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
    }
}

# ---------------------------------------------------------------------
# de-de sheet: row 3 is the "b.md" row.
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $newStatus
$wsDeDe.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("E3").Value = "2016-03-22 18:33:50"

foreach ($hl in $wsDeDe.Hyperlinks) {
    if ($hl.Range.Address() -eq "`$D`$3") {
        $hl.TextToDisplay = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
    }
}
